$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 71, shifting existing rows 71-113 down to 72-114
$ws.Rows.Item(71).Insert()

# Populate the newly inserted row 71 with the new data record
$ws.Cells.Item(71, 1).Value = 8
$ws.Cells.Item(71, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(71, 3).Value = "Coquimbo"
$ws.Cells.Item(71, 4).Value = 44438
$ws.Cells.Item(71, 5).Value = 4
$ws.Cells.Item(71, 6).Value = 100112031
$ws.Cells.Item(71, 7).Value = "Poroto verde"
$ws.Cells.Item(71, 8).Value = "Magnum"
$ws.Cells.Item(71, 9).Value = "Primera"
$ws.Cells.Item(71, 10).Value = 600
$ws.Cells.Item(71, 11).Value = 29000
$ws.Cells.Item(71, 12).Value = 30000
$ws.Cells.Item(71, 13).Value = 29500
$ws.Cells.Item(71, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(71, 15).Value = "Perú"
$ws.Cells.Item(71, 16).Value = 1180
$ws.Cells.Item(71, 17).Value = 25
$ws.Cells.Item(71, 18).Value = "Hortaliza"
